# The deck currently applies the "Integral" colour theme (theme2.xml, the
# theme linked from the one slide master) to every slide; the notes master
# separately links theme1.xml, which still holds stock "Office Theme"
# colours. The edit swaps the two: the slides' theme becomes the default
# "Office" colour set while the (effectively unused) notes-master theme
# ends up holding what used to be the slide theme's colours.
#
# PowerPoint's COM object model only exposes read/write access to theme
# colours through ThemeColorScheme.Colors(index).RGB (there is no supported
# way to rename a theme or a colour scheme from VBA/COM), so this script
# reassigns the twelve theme colour slots on the presentation's slide
# master to the "Office" theme's RGB values, in the standard
# dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink order.

$p = $ppt.ActivePresentation

function ToComRGB($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

$officeColors = @(
    (ToComRGB 0x00 0x00 0x00), # 1  dk1
    (ToComRGB 0xFF 0xFF 0xFF), # 2  lt1
    (ToComRGB 0x44 0x54 0x6A), # 3  dk2
    (ToComRGB 0xE7 0xE6 0xE6), # 4  lt2
    (ToComRGB 0x5B 0x9B 0xD5), # 5  accent1
    (ToComRGB 0xED 0x7D 0x31), # 6  accent2
    (ToComRGB 0xA5 0xA5 0xA5), # 7  accent3
    (ToComRGB 0xFF 0xC0 0x00), # 8  accent4
    (ToComRGB 0x44 0x72 0xC4), # 9  accent5
    (ToComRGB 0x70 0xAD 0x47), # 10 accent6
    (ToComRGB 0x05 0x63 0xC1), # 11 hlink
    (ToComRGB 0x95 0x4F 0x72)  # 12 folHlink
)

$tcs = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $tcs.Colors($i).RGB = $officeColors[$i - 1]
}
